# Added missing e2e scenario
# - Row 5 (sampleId "NOLANEID") is disambiguated to "NOLANEID1"
# - A brand-new row 6 is appended, duplicating row 5's scenario but with
#   sampleId "1234STDY1236" / laneId "NOLANEID2"
# - The active selection moves to D11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 6: add the missing e2e scenario ------------------------------
# Populate C6 before touching D5 so the shared-string table gets the new
# unique strings appended in the same order the workbook author's Excel did
# (1234STDY1236, then NOLANEID1, then NOLANEID2).
$ws.Range("A6").Value = "MENINGITIS"
$ws.Range("C6").Value = "1234STDY1236"

# --- Existing row 5: rename the duplicate laneId placeholder ---------------
$ws.Range("D5").Value = "NOLANEID1"

# --- Finish populating new row 6 -------------------------------------------
$ws.Range("D6").Value = "NOLANEID2"
$ws.Range("E6").Value = "IB"
$ws.Range("F6").Value = "Wellcome Sanger Institute"

# Match column F's distinct (Menlo) cell style used by the other data rows.
$ws.Range("F5").Copy()
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Move the saved selection to D11 ---------------------------------------
$ws.Range("D11").Select() | Out-Null
